$wb = $excel.ActiveWorkbook

# Remember which sheet is currently active so we can restore it afterward.
$originalActiveSheet = $wb.ActiveSheet.Name

$ws = $wb.Worksheets.Item("SaleSearchTransactions")

# Add "Y" (Execute flag) to column D for all data rows (2-12, 14-18); row 13 already had it.
foreach ($r in 2..12) {
    $ws.Range("D$r").Value = "Y"
}
foreach ($r in 14..18) {
    $ws.Range("D$r").Value = "Y"
}

# Update the selection on the sheet to match the edited range.
$ws.Activate()
$ws.Range("D2:D18").Select()

# Restore the originally active sheet/tab.
$wb.Worksheets.Item($originalActiveSheet).Activate()
